$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns B..G): TB, d2S, K, IP, Win, sum
$data = @{
    2  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    3  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433)
    4  = @(0.04172184405617529, 0.3048912486333797, 3.223369029078222, 13.86384647080068, 0, 17.43382859256846)
    5  = @(0.6545652718822623, 1.626987699542094, 189.6080260415259, 13.86384647080068, 1, 205.753425483751)
    6  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 1, 21.98653043760045)
    7  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 1, 19.48425592650926)
    8  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    9  = @(0.6545652718822623, 1.626987699542094, 189.6080260415259, 13.86384647080068, 0, 205.753425483751)
    10 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 4.327115817150455)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
